$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - sheet1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1149
$ws1.Range("F8").Value = 268
$ws1.Range("F15").Value = 13010
$ws1.Range("F16").Value = 155
$ws1.Range("F18").Value = 5351
$ws1.Range("F19").Value = 5545

# Sheet "全部类型" (All Types) - sheet4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1149
$ws4.Range("F9").Value = 268
$ws4.Range("F17").Value = 13010
$ws4.Range("F18").Value = 155
$ws4.Range("F21").Value = 5351
$ws4.Range("F22").Value = 5545
